# Intro to R.pptx - refresh the auto "last saved" date/time fields
# (the deck was simply re-uploaded later; PowerPoint recalculates the
# datetimeFigureOut header/footer fields, and the title slide repeats
# that date as plain text) from 2023-11-08 to 2024-05-21.

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16

function Set-DatePlaceholderText {
    param(
        $shapes,
        [string]$newText
    )
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = -1
        try { $phType = $shp.PlaceholderFormat.Type } catch {}
        if ($phType -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# 1) Slide master date placeholder: 08/11/2023 -> 21/05/2024 (en-GB)
Set-DatePlaceholderText -shapes $p.SlideMaster.Shapes -newText "21/05/2024"

# 2) Every slide layout's date placeholder: 08/11/2023 -> 21/05/2024 (en-GB)
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText -shapes $layout.Shapes -newText "21/05/2024"
}

# 3) Notes master date placeholder: 11/8/2023 -> 5/21/2024 (en-US)
Set-DatePlaceholderText -shapes $p.NotesMaster.Shapes -newText "5/21/2024"

# 4) Title slide subtitle literally spells the date out: "8 November, 2023"
#    becomes "21 May, 2024" - edit just the changed words/number so the
#    surrounding runs (language tags, spacing, punctuation) stay intact.
$titleSlide = $p.Slides.Item(1)
$dateShape = $titleSlide.Shapes.Item(2)
$tr = $dateShape.TextFrame.TextRange

$t = $tr.Text
$dayIdx = $t.IndexOf("8 November, 2023")
$tr.Characters($dayIdx + 1, 1).Text = "21"

$t = $tr.Text
$monthIdx = $t.IndexOf("November")
$tr.Characters($monthIdx + 1, 8).Text = "May"

$t = $tr.Text
$yearIdx = $t.LastIndexOf("23")
$tr.Characters($yearIdx + 1, 2).Text = "24"
